$p = $ppt.ActivePresentation

# Slide 1 title: "Example" " " "numbering" " " "MWE" (5 runs)
#   -> "Example " "numbering " "MWE" (3 runs)
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 8).Text = "Example "
$tr1.Characters(9, 10).Text = "numbering "

# Slide 2 title: "A" " " "second" " " "slide" (5 runs)
#   -> "A " "second " "slide" (3 runs)
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 2).Text = "A "
$tr2.Characters(3, 7).Text = "second "
